$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 10.79677
$ws.Range("H2").Value = 597.8260182199999
$ws.Range("M2").Value = 0.6362839830655896
$ws.Range("N2").Value = 63.23605810606893
$ws.Range("G3").Value = 11.52918716
$ws.Range("H3").Value = 1055.53156246
$ws.Range("M3").Value = 0.9718996719642737
$ws.Range("N3").Value = 192.4284507998276
$ws.Range("G4").Value = 3.78076368
$ws.Range("H4").Value = 117.94362866
$ws.Range("M4").Value = 0.489009486348819
$ws.Range("N4").Value = 28.21880362596037
$ws.Range("G5").Value = 4.36648536
$ws.Range("H5").Value = 215.55365134
$ws.Range("M5").Value = 0.5599332302859219
$ws.Range("N5").Value = 55.77066125852885
$ws.Range("G6").Value = 1.18651352
$ws.Range("H6").Value = 19.7299739
$ws.Range("M6").Value = 0.3079555673153038
$ws.Range("N6").Value = 9.327538351386101
$ws.Range("G7").Value = 1.55814054
$ws.Range("H7").Value = 40.88412344
$ws.Range("M7").Value = 0.305981178834762
$ws.Range("N7").Value = 15.71726688320439
$ws.Range("G8").Value = 0.56608654
$ws.Range("H8").Value = 5.680898699999999
$ws.Range("M8").Value = 0.1828766325731941
$ws.Range("N8").Value = 3.214605881882596
$ws.Range("G9").Value = 0.7736206600000001
$ws.Range("H9").Value = 13.7672833
$ws.Range("M9").Value = 0.204335321022576
$ws.Range("N9").Value = 7.375229099642515
$ws.Range("G10").Value = 0.32952664
$ws.Range("H10").Value = 2.6020052
$ws.Range("M10").Value = 0.1159635243767719
$ws.Range("N10").Value = 1.498364126726572
$ws.Range("G11").Value = 0.44701162
$ws.Range("H11").Value = 6.23821166
$ws.Range("M11").Value = 0.1382148198494925
$ws.Range("N11").Value = 3.731312673529347
$ws.Range("G12").Value = 0.1998774
$ws.Range("H12").Value = 1.25345386
$ws.Range("M12").Value = 0.09040462050834856
$ws.Range("N12").Value = 0.8774745004838026
$ws.Range("G13").Value = 0.28237826
$ws.Range("H13").Value = 3.31140634
$ws.Range("M13").Value = 0.09604014740408928
$ws.Range("N13").Value = 2.392863867332192
